# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values per row; update with newly recalculated values
$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    10 = 1
    11 = 1
    15 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
